# Apply "Natmi following Dr Hou advice" update:
# - Existing rows 2-5 get updated Sending/Ligand/Receptor/Target cluster labels and recalculated
#   expression statistics.
# - Four new rows (6-9) are appended, mirroring rows 2-5 but for the "M2" sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "S100a9"
$ws.Range("C2").Value = "Tlr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 2058.152994666666
$ws.Range("H2").Value = 6174.458984
$ws.Range("I2").Value = 0.9825919405339512
$ws.Range("J2").Value = 0.982591940533951
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.93036733333333
$ws.Range("N2").Value = 44.791102
$ws.Range("O2").Value = 0.252612808865421
$ws.Range("P2").Value = 0.252612808865421
$ws.Range("Q2").Value = 30728.98023857337
$ws.Range("R2").Value = 276560.8221471604
$ws.Range("S2").Value = 0.2482153100668061
$ws.Range("T2").Value = 0.2482153100668061

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "S100a9"
$ws.Range("C3").Value = "Tlr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 2058.152994666666
$ws.Range("H3").Value = 6174.458984
$ws.Range("I3").Value = 0.9825919405339512
$ws.Range("J3").Value = 0.982591940533951
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.28486166666667
$ws.Range("N3").Value = 60.854585
$ws.Range("O3").Value = 0.3432076230048887
$ws.Range("P3").Value = 0.3432076230048887
$ws.Range("Q3").Value = 41749.34878564907
$ws.Range("R3").Value = 375744.1390708417
$ws.Range("S3").Value = 0.3372330442944184
$ws.Range("T3").Value = 0.3372330442944183

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "S100a9"
$ws.Range("C4").Value = "Tlr4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 2058.152994666666
$ws.Range("H4").Value = 6174.458984
$ws.Range("I4").Value = 0.9825919405339512
$ws.Range("J4").Value = 0.982591940533951
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.48898366666667
$ws.Range("N4").Value = 58.466951
$ws.Range("O4").Value = 0.3297418473407271
$ws.Range("P4").Value = 0.3297418473407271
$ws.Range("Q4").Value = 40111.31009655975
$ws.Range("R4").Value = 361001.7908690378
$ws.Range("S4").Value = 0.3240016816537749
$ws.Range("T4").Value = 0.3240016816537749

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "S100a9"
$ws.Range("C5").Value = "Tlr4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 2058.152994666666
$ws.Range("H5").Value = 6174.458984
$ws.Range("I5").Value = 0.9825919405339512
$ws.Range("J5").Value = 0.982591940533951
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.399549333333334
$ws.Range("N5").Value = 13.198648
$ws.Range("O5").Value = 0.07443772078896321
$ws.Range("P5").Value = 0.07443772078896321
$ws.Range("Q5").Value = 9054.945635583737
$ws.Range("R5").Value = 81494.51072025363
$ws.Range("S5").Value = 0.0731419045189518
$ws.Range("T5").Value = 0.07314190451895179

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "S100a9"
$ws.Range("C6").Value = "Tlr4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 36.46320333333333
$ws.Range("H6").Value = 109.38961
$ws.Range("I6").Value = 0.01740805946604894
$ws.Range("J6").Value = 0.01740805946604894
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.93036733333333
$ws.Range("N6").Value = 44.791102
$ws.Range("O6").Value = 0.252612808865421
$ws.Range("P6").Value = 0.252612808865421
$ws.Range("Q6").Value = 544.409019916691
$ws.Range("R6").Value = 4899.681179250219
$ws.Range("S6").Value = 0.004397498798614903
$ws.Range("T6").Value = 0.004397498798614902

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "S100a9"
$ws.Range("C7").Value = "Tlr4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 36.46320333333333
$ws.Range("H7").Value = 109.38961
$ws.Range("I7").Value = 0.01740805946604894
$ws.Range("J7").Value = 0.01740805946604894
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.28486166666667
$ws.Range("N7").Value = 60.854585
$ws.Range("O7").Value = 0.3432076230048887
$ws.Range("P7").Value = 0.3432076230048887
$ws.Range("Q7").Value = 739.6510355402056
$ws.Range("R7").Value = 6656.85931986185
$ws.Range("S7").Value = 0.005974578710470409
$ws.Range("T7").Value = 0.005974578710470408

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "S100a9"
$ws.Range("C8").Value = "Tlr4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 36.46320333333333
$ws.Range("H8").Value = 109.38961
$ws.Range("I8").Value = 0.01740805946604894
$ws.Range("J8").Value = 0.01740805946604894
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 19.48898366666667
$ws.Range("N8").Value = 58.466951
$ws.Range("O8").Value = 0.3297418473407271
$ws.Range("P8").Value = 0.3297418473407271
$ws.Range("Q8").Value = 710.6307741976789
$ws.Range("R8").Value = 6395.676967779111
$ws.Range("S8").Value = 0.005740165686952208
$ws.Range("T8").Value = 0.005740165686952207

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "S100a9"
$ws.Range("C9").Value = "Tlr4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 36.46320333333333
$ws.Range("H9").Value = 109.38961
$ws.Range("I9").Value = 0.01740805946604894
$ws.Range("J9").Value = 0.01740805946604894
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.399549333333334
$ws.Range("N9").Value = 13.198648
$ws.Range("O9").Value = 0.07443772078896321
$ws.Range("P9").Value = 0.07443772078896321
$ws.Range("Q9").Value = 160.4216619163645
$ws.Range("R9").Value = 1443.79495724728
$ws.Range("S9").Value = 0.001295816270011419
$ws.Range("T9").Value = 0.001295816270011419
